$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (content unchanged; style was already applied in the source file)
$ws.Range("B1").Value = "Win/Loss/Draw"
$ws.Range("C1").Value = "Role"

# Observation rows: an id column plus the win/loss outcome, and (only on
# the last row) the role, which is missing ("nan") for the rest.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "win"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "win"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "win"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "loss"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "loss"
$ws.Range("C6").Value = "nan"

# Apply the same formatting used on the header cells (bold font, thin box
# border, centered/top aligned) to the id column by copying the format
# from B1, matching the existing style (rather than synthesizing a new one).
$ws.Range("B1").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
